# Add "Collective Handle Position" (H_COLL_LEV) + SBALL_ANI columns to the
# Tabelle2 sim-var matrix, inserting two new columns just before the
# END_OF_COL marker column (which, together with the Title helper column,
# shifts two columns to the right).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# The existing sheet layout (before the edit) is:
#   ... FF=THOOK | FG=END_OF_COL | FH=Title (=Tabelle1!AV<row>) ...
# Insert two blank columns at FG:FH - this pushes END_OF_COL -> FI and
# Title -> FJ, and inherits formatting/column width from the column to the
# left (FF), exactly like Excel's native "Insert Columns" command.
$ws.Columns("FG:FH").Insert()

# New header row (row 1) labels for the two inserted columns.
$ws.Range("FG1").Value = "SBALL_ANI"
$ws.Range("FH1").Value = "H_COLL_LEV"

# Populate the data rows (2-41).
#  - existing H_TRIM column (EU) switches from the "|" placeholder text to
#    a numeric flag (0, except row 41 which is 1)
#  - new SBALL_ANI column (FG) gets the "|" placeholder text like its
#    neighboring flag columns
#  - new H_COLL_LEV column (FH) gets a numeric flag (0, except row 41
#    which is 1, mirroring H_TRIM)
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 151).Value = 0   # EU: H_TRIM
    $ws.Cells.Item($r, 163).Value = "|" # FG: SBALL_ANI
    $ws.Cells.Item($r, 164).Value = 0   # FH: H_COLL_LEV
}
$ws.Cells.Item(41, 151).Value = 1  # EU41: H_TRIM
$ws.Cells.Item(41, 164).Value = 1  # FH41: H_COLL_LEV

# Restore the final selection left in the file.
$null = $ws.Range("EU26").Select()
